$wb = $excel.ActiveWorkbook
$wsCauHoi = $wb.Worksheets.Item("CauHoi")
$wsCauTraLoi = $wb.Worksheets.Item("CauTraLoi")

# ---------------------------------------------------------------------------
# Sheet "CauHoi" (questions)
# ---------------------------------------------------------------------------

# Row 2 - existing question gets replaced with a new one, DoKho changes 1 -> 2
$wsCauHoi.Range("A2").Value = "Trong quy trình phát triển phần mềm, giai đoạn nào tập trung vào việc xác định các chức năng và ràng buộc của hệ thống?"
$wsCauHoi.Range("A2").Style = "Normal"
$wsCauHoi.Range("C2").Value = 2

# Row 3 - brand new question (MaCauHoi = 4)
$wsCauHoi.Range("A3").Value = "Trong mô hình thác nước (Waterfall), bước nào phải được hoàn thành trước khi chuyển sang bước tiếp theo?"
$wsCauHoi.Range("A3").Style = "Normal"
$wsCauHoi.Range("B3").Value = 3
$wsCauHoi.Range("C3").Value = 2
$wsCauHoi.Range("D3").Value = "Trắc nghiệm"

# Row 4 - brand new question (MaCauHoi = 5)
$wsCauHoi.Range("A4").Value = "Công cụ nào thường được sử dụng để quản lý phiên bản mã nguồn trong các dự án phần mềm?"
$wsCauHoi.Range("A4").Style = "Normal"
$wsCauHoi.Range("B4").Value = 3
$wsCauHoi.Range("C4").Value = 3
$wsCauHoi.Range("D4").Value = "Trắc nghiệm"

# ---------------------------------------------------------------------------
# Sheet "CauTraLoi" (answers)
# ---------------------------------------------------------------------------

# Existing rows 2-5 belong to MaCauHoi = 3, the answer texts are updated to
# match the new question 2 on CauHoi, and the correct-answer flag moves.
$wsCauTraLoi.Range("B2").Value = "Thiết kế hệ thống (System Design)"
$wsCauTraLoi.Range("C2").Value = 0

$wsCauTraLoi.Range("B3").Value = "Phân tích yêu cầu (Requirement Analysis)"
$wsCauTraLoi.Range("C3").Value = 1

$wsCauTraLoi.Range("B4").Value = "Lập trình (Programming)"
$wsCauTraLoi.Range("C4").Value = 0

$wsCauTraLoi.Range("B5").Value = "Kiểm thử (Testing)"
$wsCauTraLoi.Range("C5").Value = 0

# New answers for question MaCauHoi = 4 (Waterfall)
$wsCauTraLoi.Range("A6").Value = 4
$wsCauTraLoi.Range("B6").Value = "Lập kế hoạch (Planning)"
$wsCauTraLoi.Range("C6").Value = 0

$wsCauTraLoi.Range("A7").Value = 4
$wsCauTraLoi.Range("B7").Value = "Yêu cầu (Requirement)"
$wsCauTraLoi.Range("C7").Value = 0

$wsCauTraLoi.Range("A8").Value = 4
$wsCauTraLoi.Range("B8").Value = "Triển khai (Deployment)"
$wsCauTraLoi.Range("C8").Value = 0

$wsCauTraLoi.Range("A9").Value = 4
$wsCauTraLoi.Range("B9").Value = "Tất cả các bước trên"
$wsCauTraLoi.Range("C9").Value = 1

# New answers for question MaCauHoi = 5 (version control tool)
$wsCauTraLoi.Range("A10").Value = 5
$wsCauTraLoi.Range("B10").Value = "Docker"
$wsCauTraLoi.Range("C10").Value = 0

$wsCauTraLoi.Range("A11").Value = 5
$wsCauTraLoi.Range("B11").Value = "Git"
$wsCauTraLoi.Range("C11").Value = 1

$wsCauTraLoi.Range("A12").Value = 5
$wsCauTraLoi.Range("B12").Value = "Jenkins"
$wsCauTraLoi.Range("C12").Value = 0

$wsCauTraLoi.Range("A13").Value = 5
$wsCauTraLoi.Range("B13").Value = "Jira"
$wsCauTraLoi.Range("C13").Value = 0

# ---------------------------------------------------------------------------
# Selections: CauHoi ends with A4 selected, CauTraLoi (the active tab) ends
# with N23 selected - select on CauHoi first, then re-activate CauTraLoi so
# it remains the active sheet (matches original tabSelected/activeTab state).
# ---------------------------------------------------------------------------
$wsCauHoi.Range("A4").Select() | Out-Null
$wsCauTraLoi.Activate() | Out-Null
$wsCauTraLoi.Range("N23").Select() | Out-Null
